$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Center the "Periodo Mora" column for the existing worker rows so the new rows
# (copied from these) pick up the same alignment.
$ws.Range("E16:E19").HorizontalAlignment = -4108

# Insert two new rows right after the existing data rows (after row 19, before the
# old row 20) so the footer/legend rows shift down to make room for a new period.
$ws.Rows("20:21").Insert()

# Duplicate the formatting (borders, fills, fonts, number formats) of the previous
# period's two worker rows (18:19, period 2508) onto the freshly inserted rows.
$ws.Range("B18:J19").Copy($ws.Range("B20:J21"))

# The copied rows still say period "2508" - update them to the new period "2509".
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"

# Update the account summary: total overdue value and count of periods now that a
# third period (2509) has been added.
$ws.Range("E11").Value = 322350
$ws.Range("F13").Value = 3
